$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 8850211.927541619

$ws.Range("B2:F7").Value = $newValue
